$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-09 Wednesday" "2023-08-10 Thursday"

Replace-Text "57×42=2394" "30×54=1620"
Replace-Text "90×96=8640" "95×72=6840"
Replace-Text "88×38=3344" "27×48=1296"
Replace-Text "21×26=546" "95×39=3705"
Replace-Text "16×40=640" "82×25=2050"

Replace-Text "92×41=3772" "31×62=1922"
Replace-Text "74×20=1480" "95×58=5510"
Replace-Text "52×62=3224" "37×65=2405"
Replace-Text "18×37=666" "77×53=4081"
Replace-Text "64×69=4416" "82×23=1886"

Replace-Text "98×12=1176" "40×72=2880"
Replace-Text "17×17=289" "80×57=4560"
Replace-Text "69×18=1242" "12×94=1128"
Replace-Text "98×70=6860" "64×47=3008"
Replace-Text "27×99=2673" "44×82=3608"

Replace-Text "89×78=6942" "87×71=6177"
Replace-Text "72×70=5040" "26×73=1898"
Replace-Text "29×31=899" "78×95=7410"
Replace-Text "86×43=3698" "54×91=4914"
Replace-Text "46×13=598" "31×12=372"

Replace-Text "83×80=6640" "78×57=4446"
Replace-Text "27×80=2160" "93×97=9021"
Replace-Text "18×62=1116" "61×95=5795"
Replace-Text "48×33=1584" "68×27=1836"
Replace-Text "58×27=1566" "93×90=8370"
